$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest run
$ws.Range("D2").Value = "90.354.05"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").Value = "3.144.77"
$ws.Range("E3").Value = "  +4.04%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "622.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +26.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.365"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.141.89"
$ws.Range("E10").Value = "  +3.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.750"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.199"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000242"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.12%  "
$ws.Range("E15").Value = "  +11.25%  "
$ws.Range("D16").Value = "90.347.16"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "3.728.87"
$ws.Range("E17").Value = "  +3.55%  "
$ws.Range("D18").Value = "3.180.93"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000209"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "454.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.91%  "
$ws.Range("D28").Value = "3.320.19"
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.160"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.20%  "
$ws.Range("E32").Value = "  -8.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +21.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "511.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.41%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.142"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.53%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.97%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.176"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +29.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0884"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +30.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.414"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "146.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  +12.97%  "
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.656"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.09%  "

Write-Host "Updated cryptos list"
